$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2019, 5.83, 1052, 4814, 1304, 158, 775, 368, 43, 1, 776, 681),
    @(2018, 4.82, 1062.2, 4792, 1263, 93, 858, 369, 54, 2, 659, 569),
    @(2017, 5.37, 1065.0999999999999, 4816, 1230, 129, 867, 431, 58, 1, 711, 636),
    @(2016, 6.1, 1054, 4947, 1372, 155, 831, 439, 54, 1, 838, 714),
    @(2015, 5.1100000000000003, 1060.0999999999999, 4835, 1237, 110, 699, 480, 68, 5, 713, 602),
    @(2014, 3.73, 1076.2, 4670, 1180, 57, 622, 330, 48, 0, 537, 446)
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 1; $col -le $rowData.Count; $col++) {
        $ws.Cells.Item($row, $col).Value = $rowData[$col - 1]
    }
}

$ws.Range("D12").Select()
